$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the bold/border/
# centered style used by the other header cells (e.g. H1) so the new
# columns match the existing header formatting exactly.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate the new I and J data columns for rows 2-73.
$ws.Cells.Item(2, 9).Value = 3
$ws.Cells.Item(2, 10).Value = 3
$ws.Cells.Item(3, 9).Value = 6
$ws.Cells.Item(3, 10).Value = 6
$ws.Cells.Item(4, 9).Value = 6
$ws.Cells.Item(4, 10).Value = 7
$ws.Cells.Item(5, 9).Value = 8
$ws.Cells.Item(5, 10).Value = 9
$ws.Cells.Item(6, 9).Value = 2
$ws.Cells.Item(6, 10).Value = 2
$ws.Cells.Item(7, 9).Value = 5
$ws.Cells.Item(7, 10).Value = 6
$ws.Cells.Item(8, 9).Value = 7
$ws.Cells.Item(8, 10).Value = 7
$ws.Cells.Item(9, 9).Value = 8
$ws.Cells.Item(9, 10).Value = 8
$ws.Cells.Item(10, 9).Value = 7
$ws.Cells.Item(10, 10).Value = 7
$ws.Cells.Item(11, 9).Value = 8
$ws.Cells.Item(11, 10).Value = 8
$ws.Cells.Item(12, 9).Value = 7
$ws.Cells.Item(12, 10).Value = 7
$ws.Cells.Item(13, 9).Value = 8
$ws.Cells.Item(13, 10).Value = 8
$ws.Cells.Item(14, 9).Value = 8
$ws.Cells.Item(14, 10).Value = 8
$ws.Cells.Item(15, 9).Value = 7
$ws.Cells.Item(15, 10).Value = 7
$ws.Cells.Item(16, 9).Value = 5
$ws.Cells.Item(16, 10).Value = 5
$ws.Cells.Item(17, 9).Value = 8
$ws.Cells.Item(17, 10).Value = 8
$ws.Cells.Item(18, 9).Value = 8
$ws.Cells.Item(18, 10).Value = 8
$ws.Cells.Item(19, 9).Value = 7
$ws.Cells.Item(19, 10).Value = 8
$ws.Cells.Item(20, 9).Value = 6
$ws.Cells.Item(20, 10).Value = 7
$ws.Cells.Item(21, 9).Value = 6
$ws.Cells.Item(21, 10).Value = 6
$ws.Cells.Item(22, 9).Value = 6
$ws.Cells.Item(22, 10).Value = 6
$ws.Cells.Item(23, 9).Value = 6
$ws.Cells.Item(23, 10).Value = 6
$ws.Cells.Item(24, 9).Value = 8
$ws.Cells.Item(24, 10).Value = 8
$ws.Cells.Item(25, 9).Value = 9
$ws.Cells.Item(25, 10).Value = 9
$ws.Cells.Item(26, 9).Value = 10
$ws.Cells.Item(26, 10).Value = 10
$ws.Cells.Item(27, 9).Value = 7
$ws.Cells.Item(27, 10).Value = 7
$ws.Cells.Item(28, 9).Value = 6
$ws.Cells.Item(28, 10).Value = 6
$ws.Cells.Item(29, 9).Value = 7
$ws.Cells.Item(29, 10).Value = 7
$ws.Cells.Item(30, 9).Value = 6
$ws.Cells.Item(30, 10).Value = 6
$ws.Cells.Item(31, 9).Value = 8
$ws.Cells.Item(31, 10).Value = 8
$ws.Cells.Item(32, 9).Value = 10
$ws.Cells.Item(32, 10).Value = 10
$ws.Cells.Item(33, 9).Value = 6
$ws.Cells.Item(33, 10).Value = 6
$ws.Cells.Item(34, 9).Value = 5
$ws.Cells.Item(34, 10).Value = 5
$ws.Cells.Item(35, 9).Value = 7
$ws.Cells.Item(35, 10).Value = 8
$ws.Cells.Item(36, 9).Value = 8
$ws.Cells.Item(36, 10).Value = 8
$ws.Cells.Item(37, 9).Value = 9
$ws.Cells.Item(37, 10).Value = 9
$ws.Cells.Item(38, 9).Value = 6
$ws.Cells.Item(38, 10).Value = 6
$ws.Cells.Item(39, 9).Value = 7
$ws.Cells.Item(39, 10).Value = 7
$ws.Cells.Item(40, 9).Value = 6
$ws.Cells.Item(40, 10).Value = 7
$ws.Cells.Item(41, 9).Value = 7
$ws.Cells.Item(41, 10).Value = 8
$ws.Cells.Item(42, 9).Value = 8
$ws.Cells.Item(42, 10).Value = 8
$ws.Cells.Item(43, 9).Value = 6
$ws.Cells.Item(43, 10).Value = 6
$ws.Cells.Item(44, 9).Value = 11
$ws.Cells.Item(44, 10).Value = 11
$ws.Cells.Item(45, 9).Value = 7
$ws.Cells.Item(45, 10).Value = 7
$ws.Cells.Item(46, 9).Value = 5
$ws.Cells.Item(46, 10).Value = 6
$ws.Cells.Item(47, 9).Value = 8
$ws.Cells.Item(47, 10).Value = 8
$ws.Cells.Item(48, 9).Value = 6
$ws.Cells.Item(48, 10).Value = 6
$ws.Cells.Item(49, 9).Value = 8
$ws.Cells.Item(49, 10).Value = 8
$ws.Cells.Item(50, 9).Value = 11
$ws.Cells.Item(50, 10).Value = 11
$ws.Cells.Item(51, 9).Value = 6
$ws.Cells.Item(51, 10).Value = 7
$ws.Cells.Item(52, 9).Value = 9
$ws.Cells.Item(52, 10).Value = 9
$ws.Cells.Item(53, 9).Value = 8
$ws.Cells.Item(53, 10).Value = 8
$ws.Cells.Item(54, 9).Value = 7
$ws.Cells.Item(54, 10).Value = 7
$ws.Cells.Item(55, 9).Value = 8
$ws.Cells.Item(55, 10).Value = 8
$ws.Cells.Item(56, 9).Value = 8
$ws.Cells.Item(56, 10).Value = 8
$ws.Cells.Item(57, 9).Value = 8
$ws.Cells.Item(57, 10).Value = 8
$ws.Cells.Item(58, 9).Value = 8
$ws.Cells.Item(58, 10).Value = 9
$ws.Cells.Item(59, 9).Value = 7
$ws.Cells.Item(59, 10).Value = 8
$ws.Cells.Item(60, 9).Value = 7
$ws.Cells.Item(60, 10).Value = 8
$ws.Cells.Item(61, 9).Value = 5
$ws.Cells.Item(61, 10).Value = 5
$ws.Cells.Item(62, 9).Value = 6
$ws.Cells.Item(62, 10).Value = 6
$ws.Cells.Item(63, 9).Value = 8
$ws.Cells.Item(63, 10).Value = 8
$ws.Cells.Item(64, 9).Value = 9
$ws.Cells.Item(64, 10).Value = 9
$ws.Cells.Item(65, 9).Value = 5
$ws.Cells.Item(65, 10).Value = 5
$ws.Cells.Item(66, 9).Value = 8
$ws.Cells.Item(66, 10).Value = 9
$ws.Cells.Item(67, 9).Value = 7
$ws.Cells.Item(67, 10).Value = 7
$ws.Cells.Item(68, 9).Value = 5
$ws.Cells.Item(68, 10).Value = 5
$ws.Cells.Item(69, 9).Value = 5
$ws.Cells.Item(69, 10).Value = 6
$ws.Cells.Item(70, 9).Value = 5
$ws.Cells.Item(70, 10).Value = 5
$ws.Cells.Item(71, 9).Value = 7
$ws.Cells.Item(71, 10).Value = 7
$ws.Cells.Item(72, 9).Value = 4
$ws.Cells.Item(72, 10).Value = 4
$ws.Cells.Item(73, 9).Value = 5
$ws.Cells.Item(73, 10).Value = 5

$ws.Application.CutCopyMode = $false
